# Update odds values on row 2 (Sheet1) to reflect refreshed FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value  = 2.15
$ws.Range("I2").Value  = 3.4
$ws.Range("K2").Value  = 2.05
$ws.Range("M2").Value  = 1.07
$ws.Range("N2").Value  = 9
$ws.Range("Q2").Value  = 2.1
$ws.Range("R2").Value  = 1.7
$ws.Range("X2").Value  = 10
$ws.Range("Y2").Value  = 9.5
$ws.Range("AC2").Value = 8.5
$ws.Range("AG2").Value = 9.5
$ws.Range("AI2").Value = 12
$ws.Range("AJ2").Value = 34
$ws.Range("AU2").Value = 5
$ws.Range("AV2").Value = 19
$ws.Range("AY2").Value = 81
$ws.Range("BB2").Value = 201
